# adjust error code and message for connector test with not exist order input
#
# Summary of the semantic change applied below (derived from the target diff):
#  - The old validation message "fields false not exist in entity!" is replaced
#    everywhere by the shorter "not exist in entity!".
#  - The old order-validation messages "bad request:error order: 12321" and
#    "bad request:error order: *" (error code 106601) are replaced by the same
#    "not exist in entity!" message with error code 106107, matching the
#    fields-validation error.
#  - Two of the previously "good request" order test rows (var9 / var10, for
#    order values "Siid,12321" and "updateTime,12321") now also expect the
#    same bad-request outcome ("bad request, order input contains invalid
#    field name" / 106107 / "not exist in entity!").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 24/25: iems-connector-test-7-var1/var2 -----------------------
# message text shortened; error code (L) and everything else unchanged.
$ws.Range("M24").Value = "not exist in entity!"
$ws.Range("M25").Value = "not exist in entity!"

# --- rows 34/35: iems-connector-test-8-var9/var10 -----------------------
# these used to be "good request" rows; they now describe a bad request
# caused by an invalid field name inside the order input.
$ws.Range("B34").Value = "bad request, order input contains invalid field name"
$ws.Range("L34").Value = 106107
$ws.Range("M34").Value = "not exist in entity!"

$ws.Range("B35").Value = "bad request, order input contains invalid field name"
$ws.Range("L35").Value = 106107
$ws.Range("M35").Value = "not exist in entity!"

# --- rows 36-39: iems-connector-test-9-var1..var4 -----------------------
# error code changed from 106601 to 106107 and message unified to
# "not exist in entity!"
$ws.Range("L36").Value = 106107
$ws.Range("M36").Value = "not exist in entity!"

$ws.Range("L37").Value = 106107
$ws.Range("M37").Value = "not exist in entity!"

$ws.Range("L38").Value = 106107
$ws.Range("M38").Value = "not exist in entity!"

$ws.Range("L39").Value = 106107
$ws.Range("M39").Value = "not exist in entity!"

# --- cosmetic view-state (best effort; sheet selection / scroll) --------
$ws.Range("M31").Select()
